$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: height tweak (content unchanged) ---
$ws.Rows.Item(8).RowHeight = 30

# --- Rows 12-15: fill in the new time-log entries ---
# Copy the date cell's number format (s=2, m/d/yyyy, bordered) onto the new date cells.
$ws.Range("A6").Copy()
$ws.Range("A12:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 12
$ws.Range("A12").Value = 41222
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Kohdetiedoston ""dummy"" kirjoitus (Huffman-koodit dumppautuu nyt kohdetiedostoon). Ennen kaikkea kuitenkin refaktorointia ja koodin siistimistä."
$ws.Rows.Item(12).RowHeight = 45

# Row 13
$ws.Range("A13").Value = 41223
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "Kohdetiedoston ""oikea"" kirjoitus. Sanastoa ei vielä kirjoiteta tiedoston alkuun joten purkaminen ei onnistu."
$ws.Rows.Item(13).RowHeight = 30

# Row 14
$ws.Range("A14").Value = 41224
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "Sanaston kirjoitus tiedoston alkuun. Purkualgoritmin aloitus."
$ws.Rows.Item(14).RowHeight = 30

# Row 15
$ws.Range("A15").Value = 41225
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "Purkualgoritmin debuggausta. Purkualgoritmi mahdollisesti toimiikin jo oikein mutta Huffman-koodin pakkaus tavujen biteiksi lienee buginen. Nyt sekä pakkaus että purku toimii melkein: purussa vielä bugi joka tuottaa epämääräisiä virheitä striimin sekaan. Ongelma liittyy Huffman-koodiin joka luetaan kahdessa lohkossa levyltä (todiste: kun lohkokoko > pakatun tiedoston koko, ongelma häviää)."
$ws.Rows.Item(15).RowHeight = 105

# --- Selection moves from C12 to A15 (where the new data was entered last) ---
$ws.Range("A15").Select() | Out-Null
